$wb = $excel.ActiveWorkbook

# Update "想去人数" (number of people wanting to go) counts that changed
# on the source site between scrapes.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 468
$ws1.Range("F5").Value = 131

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 468
$ws4.Range("F5").Value = 131
